$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised GDP figures for 2020-09-30 column (I), preserving text cell type
$updates = @{
    'I2' = '21170.3'
    'I3' = '18596.5'
    'I4' = '1687.4'
    'I6' = '1430.2'
    'I7' = '207.0'
    'I8' = '506.6'
    'I9' = '1344.5'
    'I11' = '13001'
    'I12' = '90.2'
    'I13' = '3494888'
    'I16' = ' 775154'
    'I19' = '1213480'
    'I25' = '1487035'
    'I26' = '1335133'
    'I28' = '828480'
    'I29' = '683703'
    'I30' = '1993051'
    'I31' = '2329605'
    'I32' = '20.249'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
